# Auto-generated Excel COM-interop edit script
# Applies cryptos list update per commit diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextCell $ws "D2" "67.580.54"
Set-TextCell $ws "D3" "3.493.82"
$ws.Range("E3").Value = "  +5.17%  "
Set-TextCell $ws "D4" "1.00"
$ws.Range("E4").Value = "  -0.07%  "
Set-TextCell $ws "D5" "586.60"
$ws.Range("E5").Value = "  +6.24%  "
Set-TextCell $ws "D6" "185.05"
$ws.Range("E6").Value = "  +7.06%  "
Set-TextCell $ws "D7" "0.635"
$ws.Range("E7").Value = "  +2.81%  "
Set-TextCell $ws "D8" "3.492.55"
$ws.Range("E8").Value = "  +5.43%  "
$ws.Range("E9").Value = "  -0.02%  "
Set-TextCell $ws "D10" "0.175"
$ws.Range("E10").Value = "  +2.91%  "
Set-TextCell $ws "D11" "0.652"
$ws.Range("E11").Value = "  +3.90%  "
Set-TextCell $ws "D12" "56.41"
$ws.Range("E12").Value = "  +5.80%  "
Set-TextCell $ws "D13" "0.0000281"
$ws.Range("E13").Value = "  +0.15%  "
Set-TextCell $ws "D14" "9.48"
$ws.Range("E14").Value = "  +5.01%  "
Set-TextCell $ws "D15" "4.036.62"
$ws.Range("E15").Value = "  +4.83%  "
Set-TextCell $ws "D16" "18.78"
$ws.Range("E16").Value = "  +3.95%  "
Set-TextCell $ws "D17" "3.484.54"
Set-TextCell $ws "D18" "67.456.65"
$ws.Range("E18").Value = "  +5.26%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell $ws "D19" "12.19"
$ws.Range("E19").Value = "  +4.65%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws "D20" "0.118"
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("E21").Value = "  +4.08%  "
Set-TextCell $ws "D22" "491.99"
$ws.Range("E22").Value = "  +8.15%  "
Set-TextCell $ws "D23" "5.46"
$ws.Range("E23").Value = "  +7.82%  "
Set-TextCell $ws "D24" "16.93"
$ws.Range("E24").Value = "  +22.66%  "
Set-TextCell $ws "D25" "4.46"
$ws.Range("E25").Value = "  +10.19%  "
Set-TextCell $ws "D26" "90.63"
$ws.Range("E26").Value = "  +4.65%  "
Set-TextCell $ws "D27" "2.97"
$ws.Range("E28").Value = "  +4.12%  "
Set-TextCell $ws "D29" "9.18"
$ws.Range("E29").Value = "  +7.36%  "
Set-TextCell $ws "D30" "31.72"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("E31").Value = "  +10.37%  "
Set-TextCell $ws "D32" "11.77"
$ws.Range("E32").Value = "  +3.37%  "
Set-TextCell $ws "D33" "64.45"
$ws.Range("E33").Value = "  +4.78%  "
Set-TextCell $ws "D34" "596.31"
$ws.Range("E34").Value = "  +5.33%  "
$ws.Range("E35").Value = "  +5.31%  "
Set-TextCell $ws "D36" "0.150"
$ws.Range("E36").Value = "  +6.78%  "
$ws.Range("E37").Value = "  -0.10%  "
Set-TextCell $ws "D38" "36.71"
$ws.Range("E38").Value = "  +4.47%  "
Set-TextCell $ws "D39" "0.389"
$ws.Range("E39").Value = "  +6.82%  "
$ws.Range("E40").Value = "  +0.84%  "
Set-TextCell $ws "D41" "0.0₃0772"
$ws.Range("E41").Value = "  +6.20%  "
Set-TextCell $ws "D42" "3.261.32"
$ws.Range("E42").Value = "  +7.19%  "
Set-TextCell $ws "D43" "2.92"
$ws.Range("E43").Value = "  +6.35%  "
$ws.Range("E44").Value = "  +3.80%  "
$ws.Range("E45").Value = "  +4.53%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell $ws "D46" "3.25"
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell $ws "D47" "2.76"
$ws.Range("E47").Value = "  +21.90%  "
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws "D49" "3.25"
$ws.Range("E49").Value = "  +11.52%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell $ws "D50" "8.76"
$ws.Range("E50").Value = "  +7.81%  "
$ws.Range("E51").Value = "  -0.15%  "
